# Update gh-pages to output generated at 456a3b4
# F2 ("想去人数" for the 三月三国潮动漫节 event) changes 5561 -> 5574
# F5 ("想去人数" for the 恋与深空only event) changes 12 -> 15
# Both values live on the "展览" sheet and are mirrored on the
# aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5574
    $ws.Range("F5").Value = 15
}
